$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Legionowo"
$ws.Range("B2").Value = "23.65℃"
$ws.Range("C2").Value = "39/zachmurzenie"
$ws.Range("D2").Value = "clear sky"
$ws.Range("E2").Value = "1018 hPa"
$ws.Range("H2").Value = "2023-07-18 20:49:17"
